$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. The long Kyrgyz title in A1 loses its trailing newline.
#    Re-assigning the cell value makes the host recompute the shared-string
#    table: the now-unused original string is dropped and the corrected
#    text is appended as a brand-new shared string (matching the diff,
#    which moves this <si> to the end of sharedStrings.xml without
#    xml:space="preserve").
$ws.Range("A1").Value = "8.3.1.2 Экономикадагы иш менен камсыз болгон бардык калктын чакан жана орто ишканаларда иштегендердин үлүшү"

# 2. Row-height tweaks for the header row and the two data rows.
$ws.Rows.Item(1).RowHeight = 45
$ws.Rows.Item(5).RowHeight = 17.25
$ws.Rows.Item(6).RowHeight = 17.25

# 3. Clear the stray saved selection (was N10) back to the top-left cell.
$excel.Goto($ws.Range("A1"))

# 4. Add the new 2023 column (N), reusing the formatting of the existing
#    2022 column (M) so no new cell styles are created.
$ws.Range("M3:M6").Copy()
$ws.Range("N3:N6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("N4").Value = 2023
$ws.Range("N5").Value = 2.5449890821474286
$ws.Range("N6").Value = 1.4569686017619159
